$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = "nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("G4").Value = "nourhan.mahmoud@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hananragab@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G5").Value = "servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G10").Value = "aml.awwad@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, wessam.atef@med.asu.edu.eg"
$ws.Range("G15").Value = "mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
$ws.Range("G16").Value = "Mayarembaby@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G17").Value = "Mayarembaby@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
$ws.Range("G27").Value = "nourhanmohamed@med.asu.edu.eg, System, randa.rabea@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G29").Value = "nahla91@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G30").Value = "nahla91@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("G34").Value = "nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("G35").Value = "gehanadel@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G36").Value = "servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G38").Value = "Fatmaelhady@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg"
$ws.Range("G46").Value = "Alshimaa_khaled@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg"
$ws.Range("G48").Value = "Mayarembaby@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
$ws.Range("G58").Value = "nourhanmohamed@med.asu.edu.eg, System, randa.rabea@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G60").Value = "nahla91@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G61").Value = "nahla91@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("G64").Value = "Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg"
$ws.Range("G65").Value = "majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg"
$ws.Range("G66").Value = "nourhan.mahmoud@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hananragab@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G69").Value = "Fatmaelhady@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg"
$ws.Range("G71").Value = "Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg"
$ws.Range("G72").Value = "aml.awwad@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, wessam.atef@med.asu.edu.eg"
$ws.Range("G76").Value = "mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
$ws.Range("G77").Value = "Alshimaa_khaled@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg"
$ws.Range("G79").Value = "Mayarembaby@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
$ws.Range("G91").Value = "nahla91@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G92").Value = "nahla91@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("G95").Value = "servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg"
$ws.Range("G96").Value = "nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm"
$ws.Range("G97").Value = "gehanadel@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G102").Value = "aml.awwad@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G108").Value = "mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
$ws.Range("G109").Value = "Mayarembaby@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G120").Value = "nourhan.mostafa@med.asu.edu.eg, dr_ahmed_mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G122").Value = "Salma.hassan@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, eman.samir@med.asu.edu.eg"
$ws.Range("G123").Value = "nahla91@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("G126").Value = "servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg"
$ws.Range("G127").Value = "nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm"
$ws.Range("G128").Value = "gehanadel@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G131").Value = "Fatmaelhady@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg"
$ws.Range("G133").Value = "aml.awwad@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G138").Value = "mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
$ws.Range("G139").Value = "Mayarembaby@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G140").Value = "Alshimaa_khaled@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg"
$ws.Range("G151").Value = "Amr-Saeed@med.asu.edu.eg, enas.omran@med.asu.edu.eg, takyas@med.asu.edu.eg, user@user.com"
$ws.Range("G153").Value = "Salma.hassan@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, eman.samir@med.asu.edu.eg"
$ws.Range("G154").Value = "Aya_hamed@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G157").Value = "eman.tantawi@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("G158").Value = "majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg"
$ws.Range("G159").Value = "gehanadel@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G162").Value = "Fatmaelhady@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg"
$ws.Range("G164").Value = "Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg"
$ws.Range("G170").Value = "Mayarembaby@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G182").Value = "nourhan.mostafa@med.asu.edu.eg, dr_ahmed_mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G184").Value = "Aya_hamed@med.asu.edu.eg, nahla91@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G185").Value = "Aya_hamed@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
